$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename "Sheet1" -> "Data"
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Sheet1")
$wsData.Name = "Data"

$wsInstr = $wb.Worksheets.Item("Instructions")

# ---------------------------------------------------------------------
# 2) Fix the double "the the" typo in the short description cell (A1)
# ---------------------------------------------------------------------
$reqCell = $wsInstr.Range("A1")
$reqCell.Value = "Some requirements for the upload file."

# ---------------------------------------------------------------------
# 3) Rewrite the file-format bullet list inside the big instructions
#    text box (A3), keeping the rich-text run structure / fonts intact.
# ---------------------------------------------------------------------
$bigCell = $wsInstr.Range("A3")
$fullText = $bigCell.Characters().Text

$oldBullets = "- Excel`n- CSV`n"
$startIdx0 = $fullText.IndexOf($oldBullets)
$start1Based = $startIdx0 + 1

$newBullets = "- Excel: file extensions xlsx, xls and xlsm (or capitalized)`n- CSV: file extension csv"
$extraRun = " (or CSV)`n"

# Replace "- Excel`n- CSV`n" with the longer description (no trailing newline).
$run = $bigCell.Characters($start1Based, $oldBullets.Length)
$run.Text = $newBullets

# Insert a brand new run right after it for " (or CSV)\n", on its own
# font (12pt Arial, matching the surrounding body-text size).
$insertAt = $start1Based + $newBullets.Length
$insertedRun = $bigCell.Characters($insertAt, 0)
$insertedRun.Text = $extraRun

# Re-apply explicit formatting to every logical run so the engine keeps
# them as distinct <r> runs instead of collapsing them into one.
$pos = 1

$titleLen = $startIdx0  # chars before "- Excel..." belong to the bold title run
$r = $bigCell.Characters($pos, $titleLen)
$r.Font.Name = "Arial"
$r.Font.Size = 15
$r.Font.Bold = $true
$pos += $titleLen

$r = $bigCell.Characters($pos, $newBullets.Length)
$r.Font.Name = "Arial"
$r.Font.Size = 12
$r.Font.Bold = $false
$pos += $newBullets.Length

$r = $bigCell.Characters($pos, $extraRun.Length)
$r.Font.Name = "Arial"
$r.Font.Size = 12
$r.Font.Bold = $false
$pos += $extraRun.Length

$fullTextNow = $bigCell.Characters().Text
$remainingLen = $fullTextNow.Length - $pos + 1

$r = $bigCell.Characters($pos, 1)
$r.Font.Name = "Arial"
$r.Font.Size = 15
$r.Font.Bold = $false
$pos += 1

$dataHeader = "Data sheet format:`n"
$r = $bigCell.Characters($pos, $dataHeader.Length)
$r.Font.Name = "Arial"
$r.Font.Size = 15
$r.Font.Bold = $true
$pos += $dataHeader.Length

$dataBody = "- event_start : Timestamp in ISO 8601 format`n- event_value : Numeric Values`n`n"
$r = $bigCell.Characters($pos, $dataBody.Length)
$r.Font.Name = "Arial"
$r.Font.Size = 12
$r.Font.Bold = $false
$pos += $dataBody.Length

$assumptionsHeader = "Assumptions:`n"
$r = $bigCell.Characters($pos, $assumptionsHeader.Length)
$r.Font.Name = "Arial"
$r.Font.Size = 15
$r.Font.Bold = $true
$pos += $assumptionsHeader.Length

$fullTextNow = $bigCell.Characters().Text
$tailLen = $fullTextNow.Length - $pos + 1
$r = $bigCell.Characters($pos, $tailLen)
$r.Font.Name = "Arial"
$r.Font.Size = 12
$r.Font.Bold = $false

# Re-editing the text makes the engine auto-fit the row height; restore
# the original explicit row height (139.55) that the workbook shipped with.
$wsInstr.Rows.Item(3).RowHeight = 139.55

# ---------------------------------------------------------------------
# 4) Move the saved selection on both sheets back to A1
# ---------------------------------------------------------------------
[void]$wsData.Range("A1").Select()
[void]$wsInstr.Range("A1").Select()

# ---------------------------------------------------------------------
# 5) Drop the sheet protection on the Instructions sheet
# ---------------------------------------------------------------------
$wsInstr.Unprotect()

Write-Host "done"
